$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2347.7273
$ws.Range("I40").Value = 2239.2144
$ws.Range("J40").Value = 2537.625
$ws.Range("K40").Value = 2239.2144
$ws.Range("L40").Value = 2537.625
$ws.Range("M40").Value = -2064.2144
$ws.Range("N40").Value = -2887.625

$ws.Range("H76").Value = 3587297
$ws.Range("I76").Value = 4447312
$ws.Range("K76").Value = 4447312
$ws.Range("M76").Value = -4446997

$ws.Range("H79").Value = 3587297
$ws.Range("I79").Value = 4447312
$ws.Range("K79").Value = 4447312
$ws.Range("M79").Value = -4446220

$ws.Range("H95").Value = 467082.34
$ws.Range("J95").Value = 467082.34
$ws.Range("L95").Value = 467082.34
$ws.Range("N95").Value = -472574.34

$ws.Range("H106").Value = 3834125.2
$ws.Range("I106").Value = 4276182
$ws.Range("J106").Value = 2966.6667
$ws.Range("K106").Value = 4276182
$ws.Range("L106").Value = 2966.6667
$ws.Range("M106").Value = -4275551
$ws.Range("N106").Value = -4228.6667

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H138").Value = 1409.2828
$ws.Range("I138").Value = 542.125
$ws.Range("J138").Value = 2994.9429
$ws.Range("K138").Value = 1626.375
$ws.Range("L138").Value = 8984.8287
$ws.Range("M138").Value = 3513.625
$ws.Range("N138").Value = -19264.8287

$ws.Range("H141").Value = 2179.0908
$ws.Range("I141").Value = 1268.463
$ws.Range("J141").Value = 6276.9165
$ws.Range("K141").Value = 3805.389
$ws.Range("L141").Value = 18830.7495
$ws.Range("M141").Value = 1374.611
$ws.Range("N141").Value = -29190.7495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 716.625
$ws.Range("I110").Value = 609.0909
$ws.Range("K110").Value = 609.0909
$ws.Range("M110").Value = 1435.9091

$ws.Range("H132").Value = 1899.8043
$ws.Range("I132").Value = 1625.1025
$ws.Range("K132").Value = 4875.3075
$ws.Range("M132").Value = -2345.3075

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6346
$ws.Range("I86").Value = 2023.4615
$ws.Range("J86").Value = 12589.667
$ws.Range("K86").Value = 2023.4615
$ws.Range("L86").Value = 12589.667
$ws.Range("M86").Value = -900.4614999999999
$ws.Range("N86").Value = -14835.667

$ws.Range("H89").Value = 6346
$ws.Range("I89").Value = 2023.4615
$ws.Range("J89").Value = 12589.667
$ws.Range("K89").Value = 10117.3075
$ws.Range("L89").Value = 62948.335
$ws.Range("M89").Value = -4501.307499999999
$ws.Range("N89").Value = -74180.33499999999

$ws.Range("H107").Value = 834.6
$ws.Range("I107").Value = 746.6
$ws.Range("J107").Value = 1098.6
$ws.Range("K107").Value = 746.6
$ws.Range("L107").Value = 1098.6
$ws.Range("M107").Value = 1173.4
$ws.Range("N107").Value = -4938.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1656.659
$ws.Range("I31").Value = 911.95654
$ws.Range("K31").Value = 911.95654
$ws.Range("M31").Value = -616.95654

$ws.Range("H34").Value = 1656.659
$ws.Range("I34").Value = 911.95654
$ws.Range("K34").Value = 911.95654
$ws.Range("M34").Value = -709.95654

$ws.Range("H58").Value = 2154.5854
$ws.Range("I58").Value = 862.6667
$ws.Range("J58").Value = 3978.4707
$ws.Range("K58").Value = 862.6667
$ws.Range("L58").Value = 3978.4707
$ws.Range("M58").Value = -659.6667
$ws.Range("N58").Value = -4384.4707

$ws.Range("H105").Value = 666.04
$ws.Range("I105").Value = 593.56525
$ws.Range("K105").Value = 593.56525
$ws.Range("M105").Value = 1153.43475

$ws.Range("H132").Value = 2003.5358
$ws.Range("I132").Value = 1512.6595
$ws.Range("J132").Value = 4567
$ws.Range("K132").Value = 4537.9785
$ws.Range("L132").Value = 13701
$ws.Range("M132").Value = -2007.9785
$ws.Range("N132").Value = -18761

$ws.Range("H134").Value = 1884.3442
$ws.Range("I134").Value = 1118.24
$ws.Range("J134").Value = 5366.636
$ws.Range("K134").Value = 3354.72
$ws.Range("L134").Value = 16099.908
$ws.Range("M134").Value = -819.7200000000003
$ws.Range("N134").Value = -21169.908

$ws.Range("H136").Value = 2154.5854
$ws.Range("I136").Value = 862.6667
$ws.Range("J136").Value = 3978.4707
$ws.Range("K136").Value = 2588.0001
$ws.Range("L136").Value = 11935.4121
$ws.Range("M136").Value = -38.0001000000002
$ws.Range("N136").Value = -17035.4121

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83333550
$ws.Range("I2").Value = 45.57143
$ws.Range("J2").Value = 200000460
$ws.Range("K2").Value = 273.42858
$ws.Range("L2").Value = 1200002760
$ws.Range("M2").Value = -160.42858
$ws.Range("N2").Value = -1200002986

$ws.Range("H5").Value = 1079.6552
$ws.Range("I5").Value = 501.8889
$ws.Range("J5").Value = 2025.091
$ws.Range("K5").Value = 1505.6667
$ws.Range("L5").Value = 6075.272999999999
$ws.Range("M5").Value = -1393.6667
$ws.Range("N5").Value = -6299.272999999999

$ws.Range("H16").Value = 890
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 890
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2670
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -3016

$ws.Range("H20").Value = 800
$ws.Range("I20").Value = 800
$ws.Range("K20").Value = 2400
$ws.Range("M20").Value = -2173

$ws.Range("H35").Value = 1002
$ws.Range("I35").Value = 1002
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 3006
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -2718
$ws.Range("N35").ClearContents()

$ws.Range("H38").Value = 141.23529
$ws.Range("J38").Value = 163.72728
$ws.Range("L38").Value = 491.18184
$ws.Range("N38").Value = -1185.18184

$ws.Range("H135").Value = 1079.6552
$ws.Range("I135").Value = 501.8889
$ws.Range("J135").Value = 2025.091
$ws.Range("K135").Value = 4517.0001
$ws.Range("L135").Value = 18225.819
$ws.Range("M135").Value = -1982.0001
$ws.Range("N135").Value = -23295.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2881.8293
$ws.Range("I132").Value = 2734.139
$ws.Range("J132").Value = 3945.2
$ws.Range("K132").Value = 8202.417000000001
$ws.Range("L132").Value = 11835.6
$ws.Range("M132").Value = -5672.417000000001
$ws.Range("N132").Value = -16895.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7062.815
$ws.Range("I61").Value = 6849.8184
$ws.Range("J61").Value = 8000
$ws.Range("K61").Value = 6849.8184
$ws.Range("L61").Value = 8000
$ws.Range("M61").Value = -6647.8184
$ws.Range("N61").Value = -8404

$ws.Range("H93").Value = 1404.1818
$ws.Range("I93").Value = 549.8333
$ws.Range("J93").Value = 2429.4
$ws.Range("K93").Value = 549.8333
$ws.Range("L93").Value = 2429.4
$ws.Range("M93").Value = 698.1667
$ws.Range("N93").Value = -4925.4

$ws.Range("H113").Value = 7062.815
$ws.Range("I113").Value = 6849.8184
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 6849.8184
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -4679.8184
$ws.Range("N113").Value = -12340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 416.66666
$ws.Range("I113").Value = 325
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 975
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 1195
$ws.Range("N113").Value = -6140

$ws.Range("H132").Value = 10003015
$ws.Range("I132").Value = 15154834
$ws.Range("J132").Value = 2423.7058
$ws.Range("K132").Value = 45464502
$ws.Range("L132").Value = 7271.117400000001
$ws.Range("M132").Value = -45461972
$ws.Range("N132").Value = -12331.1174
